$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the updated GENESIS table id
$ws.Name = "46421-0012"

# Update the "as of" timestamp string (report regenerated a few weeks later)
$ws.Range("A372").Value = "Stand: 06.07.2020 / 15:01:49"

# Fill in the previously-unavailable May 2020 figures for the three
# metrics (flight movements / passengers / freight+mail) that had been
# placeholders ("...") pending publication.
$ws.Range("R334").Value = 1
$ws.Range("S334").Value = 2
$ws.Range("T334").Value = "-"
$ws.Range("R335").Value = 25
$ws.Range("S335").Value = 54
$ws.Range("T335").Value = "-"
$ws.Range("R336").Value = 56
$ws.Range("S336").Value = 112
$ws.Range("T336").Value = 2
$ws.Range("R337").Value = 7
$ws.Range("S337").Value = 44
$ws.Range("T337").Value = 0
$ws.Range("R338").Value = 34
$ws.Range("S338").Value = 270
$ws.Range("T338").Value = 83
$ws.Range("R339").Value = 1272
$ws.Range("S339").Value = 89706
$ws.Range("T339").Value = 525
$ws.Range("R340").Value = 192
$ws.Range("S340").Value = 10780
$ws.Range("T340").Value = 2011
$ws.Range("R341").Value = 1965
$ws.Range("S341").Value = 26099
$ws.Range("T341").Value = 69740
$ws.Range("R342").Value = 3552
$ws.Range("S342").Value = 127067
$ws.Range("T342").Value = 72361
$ws.Range("R343").Value = 1
$ws.Range("S343").Value = 2
$ws.Range("T343").Value = "-"
$ws.Range("R344").Value = 12
$ws.Range("S344").Value = 26
$ws.Range("T344").Value = "-"
$ws.Range("R345").Value = 41
$ws.Range("S345").Value = 72
$ws.Range("T345").Value = "-"
$ws.Range("R346").Value = 6
$ws.Range("S346").Value = 45
$ws.Range("T346").Value = "-"
$ws.Range("R347").Value = 36
$ws.Range("S347").Value = 211
$ws.Range("T347").Value = 14
$ws.Range("R348").Value = 1271
$ws.Range("S348").Value = 83480
$ws.Range("T348").Value = 514
$ws.Range("R349").Value = 224
$ws.Range("S349").Value = 19504
$ws.Range("T349").Value = 2133
$ws.Range("R350").Value = 2142
$ws.Range("S350").Value = 41409
$ws.Range("T350").Value = 80729
$ws.Range("R351").Value = 3733
$ws.Range("S351").Value = 144749
$ws.Range("T351").Value = 83390
$ws.Range("R352").Value = 2038
$ws.Range("S352").Value = 122
$ws.Range("T352").Value = "-"
$ws.Range("R353").Value = 1027
$ws.Range("S353").Value = 613
$ws.Range("T353").Value = 0
$ws.Range("R354").Value = 639
$ws.Range("S354").Value = 1185
$ws.Range("T354").Value = 100
$ws.Range("R355").Value = 107
$ws.Range("S355").Value = 332
$ws.Range("T355").Value = 68
$ws.Range("R356").Value = 242
$ws.Range("S356").Value = 427
$ws.Range("T356").Value = 769
$ws.Range("R357").Value = 3667
$ws.Range("S357").Value = 164844
$ws.Range("T357").Value = 6907
$ws.Range("R358").Value = 1901
$ws.Range("S358").Value = 15238
$ws.Range("T358").Value = 36908
$ws.Range("R359").Value = 3511
$ws.Range("S359").Value = 29343
$ws.Range("T359").Value = 133985
$ws.Range("R360").Value = 13132
$ws.Range("S360").Value = 212104
$ws.Range("T360").Value = 178737
$ws.Range("R361").Value = 2033
$ws.Range("S361").Value = 87
$ws.Range("T361").Value = "-"
$ws.Range("R362").Value = 992
$ws.Range("S362").Value = 531
$ws.Range("T362").Value = "-"
$ws.Range("R363").Value = 626
$ws.Range("S363").Value = 1313
$ws.Range("T363").Value = 10
$ws.Range("R364").Value = 101
$ws.Range("S364").Value = 309
$ws.Range("T364").Value = 45
$ws.Range("R365").Value = 262
$ws.Range("S365").Value = 386
$ws.Range("T365").Value = 555
$ws.Range("R366").Value = 3776
$ws.Range("S366").Value = 173720
$ws.Range("T366").Value = 5989
$ws.Range("R367").Value = 1983
$ws.Range("S367").Value = 33052
$ws.Range("T367").Value = 33535
$ws.Range("R368").Value = 3790
$ws.Range("S368").Value = 43058
$ws.Range("T368").Value = 149828
$ws.Range("R369").Value = 13563
$ws.Range("S369").Value = 252456
$ws.Range("T369").Value = 189962
